$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "14.132857333438793"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "19.078635313671487"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "16.86804957264372"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "10.352979166817128"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "10.35100322721796"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "6.5104002628833495"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "5.5784055063340965"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.308084443450454"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.941165841762583"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.6268304323492213"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.2901640142887"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.389159900756226"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.191060735235943"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.985494960541278"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.083060905109048"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.08307945869664"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "21.9791299340416"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.462258580812701"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.1144524098602515"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.72633295413688"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.876157870353966"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.686678092555963"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.0680096503521805"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.596082923040191"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "64.78256506104303"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "14.134553587646586"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "13.74467659989817"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.8351612877348416"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.528322711879049"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.34615746532614"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.7937218939612976"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.34285606471344"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.0478559782110182"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.2537511408730624"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "28.037927745968446"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.1617461901411605"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.865160109630305"
$ws.Range("D52").NumberFormat = "@"
$ws.Range("D52").Value = "2.921401860711616"
$ws.Range("D54").NumberFormat = "@"
$ws.Range("D54").Value = "20.3115961545465"
$ws.Range("D56").NumberFormat = "@"
$ws.Range("D56").Value = "4.669968356634079"
$ws.Range("D58").NumberFormat = "@"
$ws.Range("D58").Value = "5.654864944514908"
$ws.Range("D59").NumberFormat = "@"
$ws.Range("D59").Value = "4.900284055080419"
$ws.Range("D62").NumberFormat = "@"
$ws.Range("D62").Value = "4.117522396646869"
$ws.Range("D63").NumberFormat = "@"
$ws.Range("D63").Value = "3.48046501322629"
$ws.Range("D64").NumberFormat = "@"
$ws.Range("D64").Value = "12.612790636182957"
$ws.Range("D65").NumberFormat = "@"
$ws.Range("D65").Value = "12.181641524022865"
$ws.Range("D66").NumberFormat = "@"
$ws.Range("D66").Value = "16.06909372770657"
$ws.Range("D67").NumberFormat = "@"
$ws.Range("D67").Value = "83.13221356728918"
$ws.Range("D68").NumberFormat = "@"
$ws.Range("D68").Value = "96.55924705790383"
$ws.Range("D69").NumberFormat = "@"
$ws.Range("D69").Value = "2.359899243802385"
$ws.Range("D70").NumberFormat = "@"
$ws.Range("D70").Value = "3.83224625823254"
$ws.Range("D71").NumberFormat = "@"
$ws.Range("D71").Value = "2.383768508382552"
$ws.Range("D72").NumberFormat = "@"
$ws.Range("D72").Value = "3.67853795375526"
$ws.Range("D73").NumberFormat = "@"
$ws.Range("D73").Value = "12.365358993764211"
$ws.Range("D74").NumberFormat = "@"
$ws.Range("D74").Value = "4.129497221360352"
$ws.Range("D75").NumberFormat = "@"
$ws.Range("D75").Value = "1.3869765361553135"
$ws.Range("D76").NumberFormat = "@"
$ws.Range("D76").Value = "4.8377838645686735"
$ws.Range("D77").NumberFormat = "@"
$ws.Range("D77").Value = "11.91122583150462"
$ws.Range("D78").NumberFormat = "@"
$ws.Range("D78").Value = "2.9755835678819724"
$ws.Range("D79").NumberFormat = "@"
$ws.Range("D79").Value = "11.791371934329248"
$ws.Range("D80").NumberFormat = "@"
$ws.Range("D80").Value = "39.32001335034202"
$ws.Range("D81").NumberFormat = "@"
$ws.Range("D81").Value = "8.179218536621091"
$ws.Range("D82").NumberFormat = "@"
$ws.Range("D82").Value = "98.09406650013413"
$ws.Range("D83").NumberFormat = "@"
$ws.Range("D83").Value = "2.271762874101956"
$ws.Range("D84").NumberFormat = "@"
$ws.Range("D84").Value = "52.51438959676666"
$ws.Range("D85").NumberFormat = "@"
$ws.Range("D85").Value = "24.765808193381645"
$ws.Range("D86").NumberFormat = "@"
$ws.Range("D86").Value = "13.802628746663991"
$ws.Range("D87").NumberFormat = "@"
$ws.Range("D87").Value = "4.659457557810683"
$ws.Range("D88").NumberFormat = "@"
$ws.Range("D88").Value = "1.4288398482509048"
$ws.Range("D89").NumberFormat = "@"
$ws.Range("D89").Value = "2.7263747212195844"
$ws.Range("D90").NumberFormat = "@"
$ws.Range("D90").Value = "0.7005651266453613"
$ws.Range("D94").NumberFormat = "@"
$ws.Range("D94").Value = "11.798381623277185"
$ws.Range("D97").NumberFormat = "@"
$ws.Range("D97").Value = "16.43805485276187"
$ws.Range("D98").NumberFormat = "@"
$ws.Range("D98").Value = "48.10705139526472"
$ws.Range("D100").NumberFormat = "@"
$ws.Range("D100").Value = "12.20922342273179"
$ws.Range("D101").NumberFormat = "@"
$ws.Range("D101").Value = "75.73586210801837"
$ws.Range("D102").NumberFormat = "@"
$ws.Range("D102").Value = "4.005429435016678"
$ws.Range("D103").NumberFormat = "@"
$ws.Range("D103").Value = "6.24945950943237"
